# Batch claim resolution: add new sources to the Sources registry and
# patch the corresponding Claims.source_ids references that were
# previously marked UNVERIFIED.

$wb = $excel.ActiveWorkbook

# --- 1) Sources sheet: append new registry rows 123-128 (S122-S127) ---
$sources = $wb.Worksheets.Item("Sources")

$newSources = @(
    @{ Row = 123; Id = "S122"; B = "sources/articles/FeedAdditive_ROI_3to1.txt"; C = "Feed & Additive Magazine"; D = "Phytogenic Feed Additives ROI"; E = "2026-02-08" },
    @{ Row = 124; Id = "S123"; B = "sources/articles/PetFoodInd_UrbanSuburban.txt"; C = "Petfood Industry"; D = "Urban vs Suburban Purchasing Habits"; E = "2026-02-08" },
    @{ Row = 125; Id = "S124"; B = "sources/regulatory/MARA_Announcement_194_Summary.txt"; C = "MARA China"; D = "China AGP Ban Announcement 194"; E = "2020-07-01" },
    @{ Row = 126; Id = "S125"; B = "sources/reports/Sector_Deal_Multiples_2020-2024.txt"; C = "Public Financial Data"; D = "Sector Deal Multiples Assessment"; E = "2026-02-08" },
    @{ Row = 127; Id = "S126"; B = "sources/regulatory/EU_Green_Claims_Directive_Summary.txt"; C = "EU Commission"; D = "Green Claims Directive Proposal"; E = "2023-03-22" },
    @{ Row = 128; Id = "S127"; B = "sources/academic/Nutrigenomics_Review_Summary.txt"; C = "Frontiers / NIH"; D = "Nutrigenomics Review"; E = "2026-02-08" }
)

foreach ($row in $newSources) {
    $r = $row.Row
    $sources.Cells.Item($r, 1).Value = $row.Id
    $sources.Cells.Item($r, 2).Value = $row.B
    $sources.Cells.Item($r, 3).Value = $row.C
    $sources.Cells.Item($r, 4).Value = $row.D

    # Column E holds a plain "yyyy-mm-dd" label (as text, matching the
    # rest of the registry) rather than a real date serial, so force
    # text format before assigning and then drop the format override
    # again to avoid leaving a stray number-format style behind.
    $cellE = $sources.Cells.Item($r, 5)
    $cellE.NumberFormat = "@"
    $cellE.Value = $row.E
    $cellE.ClearFormats()
}

# --- 2) Claims sheet: resolve previously UNVERIFIED source_ids (col D) ---
$claims = $wb.Worksheets.Item("Claims")

$claimFixes = @(
    @{ Row = 62; SourceId = "S123" },
    @{ Row = 64; SourceId = "S122" },
    @{ Row = 73; SourceId = "S125" },
    @{ Row = 74; SourceId = "S125" },
    @{ Row = 75; SourceId = "S126" },
    @{ Row = 76; SourceId = "S127" },
    @{ Row = 77; SourceId = "S124" },
    @{ Row = 85; SourceId = "S125" },
    @{ Row = 87; SourceId = "S125" }
)

foreach ($fix in $claimFixes) {
    $claims.Cells.Item($fix.Row, 4).Value = $fix.SourceId
}
